# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new (blank) column inserted in
# front of the old "Late" column (column N), pushing "Late", the spare
# "heading" column and "Outstanding" one column to the right (N->O->P->Q).
# The two other worksheet-level changes are which sheet/cell is active:
# "Edit Repayment Schedule" becomes the active tab (cell H11 selected)
# instead of "NewLoanInput", and "Repayment schedule" ends up with G19
# selected.

$wb = $excel.ActiveWorkbook

$wsLoanInput  = $wb.Worksheets.Item("NewLoanInput")
$wsRepay      = $wb.Worksheets.Item("Repayment schedule")
$wsEditRepay  = $wb.Worksheets.Item("Edit Repayment Schedule")

# --- Repayment schedule: insert a blank column before column N ---------
$wsRepay.Columns("N:N").Insert()

# Match the new column's stored width (~11 "characters") left by the
# insert in the real workbook.
$wsRepay.Columns("N:N").ColumnWidth = 10.1666666666667

# --- Update the remembered selections on each affected sheet -----------
$wsRepay.Activate()
[void]$wsRepay.Range("G19").Select()

$wsEditRepay.Activate()
[void]$wsEditRepay.Range("H11").Select()

# "Edit Repayment Schedule" is the last sheet activated, so it becomes
# the workbook's active tab (mirrors activeTab="4" / tabSelected moving
# off of "NewLoanInput" and onto "Edit Repayment Schedule").
